# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.759.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "'2.773.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'358.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "'108.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.23%  "
$ws.Range("D7").Value = "'0.556"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").Value = "'39.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.35%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'19.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("D14").Value = "'7.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'3.211.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "'2.790.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "'0.908"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'51.602.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'3.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "'13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.82%  "
$ws.Range("D22").Value = "'0.0₃0974"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").Value = "'273.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'69.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'2.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").Value = "'26.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'0.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'0.0467"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.45%  "
$ws.Range("D32").Value = "'51.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "'33.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").Value = "'5.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("D36").Value = "'0.0833"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("E39").Value = "  -7.22%  "
$ws.Range("D40").Value = "'17.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").Value = "'125.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "'2.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("D44").Value = "'2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'21.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.23%  "
$ws.Range("D46").Value = "'2.043.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").Value = "'2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("D49").Value = "'5.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.923"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("D51").Value = "'8.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
